$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '''61.217.15'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  -2.51%  '
$ws.Range('E2').Style = 'Normal'

# Row 3
$ws.Range('D3').Value = '''3.004.00'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  -1.98%  '
$ws.Range('E3').Style = 'Normal'

# Row 4
$ws.Range('D4').Value = '''0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '''  -0.08%  '
$ws.Range('E4').Style = 'Normal'

# Row 5
$ws.Range('D5').Value = '''536.63'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  -0.12%  '
$ws.Range('E5').Style = 'Normal'

# Row 6
$ws.Range('D6').Value = '''134.88'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  +1.40%  '
$ws.Range('E6').Style = 'Normal'

# Row 7
$ws.Range('D7').Value = '''0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '''  -0.07%  '
$ws.Range('E7').Style = 'Normal'

# Row 8
$ws.Range('D8').Value = '''3.000.87'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '''  -1.88%  '
$ws.Range('E8').Style = 'Normal'

# Row 9
$ws.Range('D9').Value = '''0.496'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''  +0.73%  '
$ws.Range('E9').Style = 'Normal'

# Row 10
$ws.Range('E10').Value = '''  -3.40%  '
$ws.Range('E10').Style = 'Normal'

# Row 11
$ws.Range('D11').Value = '''6.11'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  +0.26%  '
$ws.Range('E11').Style = 'Normal'

# Row 12
$ws.Range('D12').Value = '''0.447'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  -1.12%  '
$ws.Range('E12').Style = 'Normal'

# Row 13
$ws.Range('E13').Value = '''  -0.95%  '
$ws.Range('E13').Style = 'Normal'

# Row 14
$ws.Range('D14').Value = '''34.14'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  -0.24%  '
$ws.Range('E14').Style = 'Normal'

# Row 15
$ws.Range('D15').Value = '''3.483.76'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  -1.88%  '
$ws.Range('E15').Style = 'Normal'

# Row 16
$ws.Range('E16').Value = '''  -0.55%  '
$ws.Range('E16').Style = 'Normal'

# Row 17
$ws.Range('D17').Value = '''61.152.09'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''  -2.51%  '
$ws.Range('E17').Style = 'Normal'

# Row 18
$ws.Range('D18').Value = '''3.004.25'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  -1.74%  '
$ws.Range('E18').Style = 'Normal'

# Row 19
$ws.Range('E19').Value = '''  -0.13%  '
$ws.Range('E19').Style = 'Normal'

# Row 20
$ws.Range('D20').Value = '''467.68'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  -2.95%  '
$ws.Range('E20').Style = 'Normal'

# Row 21
$ws.Range('D21').Value = '''13.27'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  -0.20%  '
$ws.Range('E21').Style = 'Normal'

# Row 22
$ws.Range('E22').Value = '''  -2.51%  '
$ws.Range('E22').Style = 'Normal'

# Row 23
$ws.Range('E23').Value = '''  -1.74%  '
$ws.Range('E23').Style = 'Normal'

# Row 24
$ws.Range('D24').Value = '''79.83'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  +0.99%  '
$ws.Range('E24').Style = 'Normal'

# Row 25
$ws.Range('D25').Value = '''12.06'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  -0.27%  '
$ws.Range('E25').Style = 'Normal'

# Row 26
$ws.Range('D26').Value = '''0.997'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  -0.22%  '
$ws.Range('E26').Style = 'Normal'

# Row 27
$ws.Range('E27').Value = '''  -0.28%  '
$ws.Range('E27').Style = 'Normal'

# Row 28
$ws.Range('D28').Value = '''7.92'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  -1.77%  '
$ws.Range('E28').Style = 'Normal'

# Row 29
$ws.Range('D29').Value = '''1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  +0.37%  '
$ws.Range('E29').Style = 'Normal'

# Row 30
$ws.Range('E30').Value = '''  +1.65%  '
$ws.Range('E30').Style = 'Normal'

# Row 31
$ws.Range('B31').Value = '''EthereumClassic'
$ws.Range('B31').Style = 'Normal'
$ws.Range('C31').Value = '''https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('C31').Style = 'Normal'
$ws.Range('D31').Value = '''25.59'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  -1.35%  '
$ws.Range('E31').Style = 'Normal'

# Row 32
$ws.Range('B32').Value = '''Mantle'
$ws.Range('B32').Style = 'Normal'
$ws.Range('C32').Value = '''https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('C32').Style = 'Normal'
$ws.Range('D32').Value = '''1.15'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  +3.78%  '
$ws.Range('E32').Style = 'Normal'

# Row 33
$ws.Range('D33').Value = '''5.51'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  +3.36%  '
$ws.Range('E33').Style = 'Normal'

# Row 34
$ws.Range('D34').Value = '''55.74'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  -2.19%  '
$ws.Range('E34').Style = 'Normal'

# Row 35
$ws.Range('D35').Value = '''2.29'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '''  -2.66%  '
$ws.Range('E35').Style = 'Normal'

# Row 36
$ws.Range('E36').Value = '''  -1.53%  '
$ws.Range('E36').Style = 'Normal'

# Row 37
$ws.Range('D37').Value = '''456.57'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  -5.73%  '
$ws.Range('E37').Style = 'Normal'

# Row 38
$ws.Range('D38').Value = '''3.204.46'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  +2.97%  '
$ws.Range('E38').Style = 'Normal'

# Row 39
$ws.Range('D39').Value = '''0.0789'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  -0.77%  '
$ws.Range('E39').Style = 'Normal'

# Row 40
$ws.Range('D40').Value = '''0.0384'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  -2.17%  '
$ws.Range('E40').Style = 'Normal'

# Row 41
$ws.Range('D41').Value = '''0.119'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  +3.09%  '
$ws.Range('E41').Style = 'Normal'

# Row 42
$ws.Range('D42').Value = '''8.15'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  +1.04%  '
$ws.Range('E42').Style = 'Normal'

# Row 43
$ws.Range('D43').Value = '''27.75'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  +13.77%  '
$ws.Range('E43').Style = 'Normal'

# Row 44
$ws.Range('D44').Value = '''2.48'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  -4.98%  '
$ws.Range('E44').Style = 'Normal'

# Row 46
$ws.Range('E46').Value = '''  -1.95%  '
$ws.Range('E46').Style = 'Normal'

# Row 47
$ws.Range('B47').Value = '''Monero'
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').Value = '''https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').Value = '''120.43'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  -0.98%  '
$ws.Range('E47').Style = 'Normal'

# Row 48
$ws.Range('B48').Value = '''Fetch.AI'
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').Value = '''https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').Value = '''2.01'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  -0.01%  '
$ws.Range('E48').Style = 'Normal'

# Row 49
$ws.Range('E49').Value = '''  +0.29%  '
$ws.Range('E49').Style = 'Normal'

# Row 50
$ws.Range('D50').Value = '''0.0₃0496'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '''  -6.93%  '
$ws.Range('E50').Style = 'Normal'

# Row 51
$ws.Range('E51').Value = '''  +7.46%  '
$ws.Range('E51').Style = 'Normal'
